$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note the "admin" marker next to the CanDeleteUser permission row (row 2)
$ws.Range("F2").Value = "admin"

# Add a new permission row for publishing events
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "HasPermissionToPublishEvent"

# Match the vertical-centered style used by the rest of the permission rows
$ws.Range("A10:B10").VerticalAlignment = -4108

# Move the active selection to the newly added cell
$ws.Range("B10").Select()
